# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data updates to Sheets/Mateus_Profits.xlsx
# (columns H-N: currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3910.4849
$ws.Range("I15").Value = 3910.4849
$ws.Range("K15").Value = 11731.4547
$ws.Range("M15").Value = -11562.4547

$ws.Range("H33").Value = 10599.322
$ws.Range("I33").Value = 11123.414
$ws.Range("J33").Value = 3000
$ws.Range("K33").Value = 11123.414
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = -10894.414
$ws.Range("N33").Value = -3458

$ws.Range("H64").Value = 9601.9375
$ws.Range("I64").Value = 6346.8335
$ws.Range("K64").Value = 6346.8335
$ws.Range("M64").Value = -6098.8335

$ws.Range("H67").Value = 9601.9375
$ws.Range("I67").Value = 6346.8335
$ws.Range("K67").Value = 6346.8335
$ws.Range("M67").Value = -5488.8335

$ws.Range("H98").Value = 3111.92
$ws.Range("I98").Value = 3111.92
$ws.Range("K98").Value = 3111.92
$ws.Range("M98").Value = -1613.92

$ws.Range("H122").Value = 3111.92
$ws.Range("I122").Value = 3111.92
$ws.Range("K122").Value = 9335.76
$ws.Range("M122").Value = -6885.76

$ws.Range("H132").Value = 1209
$ws.Range("I132").Value = 851.6957
$ws.Range("K132").Value = 2555.0871
$ws.Range("M132").Value = -25.08709999999974

$ws.Range("H135").Value = 2115.7144
$ws.Range("I135").Value = 1765.4546
$ws.Range("J135").Value = 3400
$ws.Range("K135").Value = 15889.0914
$ws.Range("L135").Value = 30600
$ws.Range("M135").Value = -13354.0914
$ws.Range("N135").Value = -35670

$ws.Range("H137").Value = 49207.07
$ws.Range("I137").Value = 49207.07
$ws.Range("K137").Value = 147621.21
$ws.Range("M137").Value = -145071.21

$ws.Range("H138").Value = 9837.799999999999
$ws.Range("J138").Value = 9847.75
$ws.Range("L138").Value = 29543.25
$ws.Range("N138").Value = -39823.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2921.589
$ws.Range("I32").Value = 2199.2856
$ws.Range("K32").Value = 2199.2856
$ws.Range("M32").Value = -1912.2856

$ws.Range("H61").Value = 7941622.5
$ws.Range("I61").Value = 11115319
$ws.Range("K61").Value = 11115319
$ws.Range("M61").Value = -11115107

$ws.Range("H63").Value = 3378.4
$ws.Range("I63").Value = 3630.8333
$ws.Range("J63").Value = 2999.75
$ws.Range("K63").Value = 3630.8333
$ws.Range("L63").Value = 2999.75
$ws.Range("M63").Value = -2944.8333
$ws.Range("N63").Value = -4371.75

$ws.Range("H66").Value = 3378.4
$ws.Range("I66").Value = 3630.8333
$ws.Range("J66").Value = 2999.75
$ws.Range("K66").Value = 18154.1665
$ws.Range("L66").Value = 14998.75
$ws.Range("M66").Value = -14722.1665
$ws.Range("N66").Value = -21862.75

$ws.Range("H74").Value = 31359.564
$ws.Range("I74").Value = 2118.6316
$ws.Range("K74").Value = 2118.6316
$ws.Range("M74").Value = -1244.6316

$ws.Range("H77").Value = 31359.564
$ws.Range("I77").Value = 2118.6316
$ws.Range("K77").Value = 10593.158
$ws.Range("M77").Value = -6225.158000000001

$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992

$ws.Range("H97").Value = 605.5454999999999
$ws.Range("I97").Value = 543.2632
$ws.Range("K97").Value = 543.2632
$ws.Range("M97").Value = -47.26319999999998

$ws.Range("H102").Value = 3423.0833
$ws.Range("J102").Value = 7000
$ws.Range("L102").Value = 7000
$ws.Range("N102").Value = -10244

$ws.Range("H132").Value = 3332.4
$ws.Range("I132").Value = 2782.7715
$ws.Range("J132").Value = 5256.1
$ws.Range("K132").Value = 8348.3145
$ws.Range("L132").Value = 15768.3
$ws.Range("M132").Value = -5818.3145
$ws.Range("N132").Value = -20828.3

$ws.Range("H136").Value = 7941622.5
$ws.Range("I136").Value = 11115319
$ws.Range("K136").Value = 33345957
$ws.Range("M136").Value = -33343407

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5191.0415
$ws.Range("I3").Value = 4607.0527
$ws.Range("J3").Value = 7410.2
$ws.Range("K3").Value = 4607.0527
$ws.Range("L3").Value = 7410.2
$ws.Range("M3").Value = -4493.0527
$ws.Range("N3").Value = -7638.2

$ws.Range("H82").Value = 43654
$ws.Range("J82").Value = 42756.668
$ws.Range("L82").Value = 42756.668
$ws.Range("N82").Value = -43522.668

$ws.Range("H85").Value = 43654
$ws.Range("J85").Value = 42756.668
$ws.Range("L85").Value = 42756.668
$ws.Range("N85").Value = -45408.668

$ws.Range("H99").Value = 3877.9
$ws.Range("I99").Value = 3398.4285
$ws.Range("J99").Value = 4996.6665
$ws.Range("K99").Value = 3398.4285
$ws.Range("L99").Value = 4996.6665
$ws.Range("M99").Value = -1900.4285
$ws.Range("N99").Value = -7992.6665

$ws.Range("H105").Value = 5161.385
$ws.Range("I105").Value = 5161.385
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 5161.385
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -3414.385
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 4504.154
$ws.Range("I107").Value = 4719.1113
$ws.Range("J107").Value = 4020.5
$ws.Range("K107").Value = 4719.1113
$ws.Range("L107").Value = 4020.5
$ws.Range("M107").Value = -2799.1113
$ws.Range("N107").Value = -7860.5

$ws.Range("H132").Value = 80000
$ws.Range("J132").Value = 80000
$ws.Range("L132").Value = 80000
$ws.Range("N132").Value = -90120

$ws.Range("H134").Value = 4069.1072
$ws.Range("I134").Value = 3932.7778
$ws.Range("K134").Value = 11798.3334
$ws.Range("M134").Value = -9263.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 49931.61
$ws.Range("J6").Value = 49927.65
$ws.Range("L6").Value = 49927.65
$ws.Range("N6").Value = -50153.65

$ws.Range("H31").Value = 4252.409
$ws.Range("I31").Value = 1111.25
$ws.Range("K31").Value = 1111.25
$ws.Range("M31").Value = -816.25

$ws.Range("H34").Value = 4252.409
$ws.Range("I34").Value = 1111.25
$ws.Range("K34").Value = 1111.25
$ws.Range("M34").Value = -909.25

$ws.Range("H58").Value = 4970.067
$ws.Range("I58").Value = 2583.4443
$ws.Range("J58").Value = 8550
$ws.Range("K58").Value = 2583.4443
$ws.Range("L58").Value = 8550
$ws.Range("M58").Value = -2380.4443
$ws.Range("N58").Value = -8956

$ws.Range("H134").Value = 6290.6177
$ws.Range("I134").Value = 3252.7646
$ws.Range("K134").Value = 9758.293799999999
$ws.Range("M134").Value = -7223.293799999999

$ws.Range("H136").Value = 4970.067
$ws.Range("I136").Value = 2583.4443
$ws.Range("J136").Value = 8550
$ws.Range("K136").Value = 7750.3329
$ws.Range("L136").Value = 25650
$ws.Range("M136").Value = -5200.3329
$ws.Range("N136").Value = -30750

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 889.6667
$ws.Range("I118").Value = 889.6667
$ws.Range("K118").Value = 2669.0001
$ws.Range("M118").Value = -1426.0001

$ws.Range("H121").Value = 11111469
$ws.Range("I121").Value = 110.5
$ws.Range("K121").Value = 331.5
$ws.Range("M121").Value = 978.5

$ws.Range("H125").Value = 9895.888999999999
$ws.Range("I125").Value = 4533.3335
$ws.Range("K125").Value = 13600.0005
$ws.Range("M125").Value = -8680.000499999998

$ws.Range("H140").Value = 2947.5386
$ws.Range("I140").Value = 3379.7778
$ws.Range("J140").Value = 1975
$ws.Range("K140").Value = 10139.3334
$ws.Range("L140").Value = 5925
$ws.Range("M140").Value = -4959.3334
$ws.Range("N140").Value = -16285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 554.8333
$ws.Range("I97").Value = 469.82352
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 469.82352
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = 26.17648000000003
$ws.Range("N97").Value = -2992

$ws.Range("H122").Value = 4468.533
$ws.Range("I122").Value = 5593.222
$ws.Range("K122").Value = 16779.666
$ws.Range("M122").Value = -14329.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1622.6666
$ws.Range("I16").Value = 1622.6666
$ws.Range("K16").Value = 1622.6666
$ws.Range("M16").Value = -1452.6666

$ws.Range("H22").Value = 2277.3333
$ws.Range("J22").Value = 2899.4
$ws.Range("L22").Value = 2899.4
$ws.Range("N22").Value = -3489.4

$ws.Range("H27").Value = 2277.3333
$ws.Range("J27").Value = 2899.4
$ws.Range("L27").Value = 2899.4
$ws.Range("N27").Value = -3113.4

$ws.Range("H100").Value = 3127995.2
$ws.Range("I100").Value = 6252473
$ws.Range("J100").Value = 3517.375
$ws.Range("K100").Value = 6252473
$ws.Range("L100").Value = 3517.375
$ws.Range("M100").Value = -6251932
$ws.Range("N100").Value = -4599.375

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 4176.6665
$ws.Range("I107").Value = 3739.4285
$ws.Range("J107").Value = 5707
$ws.Range("K107").Value = 11218.2855
$ws.Range("L107").Value = 17121
$ws.Range("M107").Value = -9298.2855
$ws.Range("N107").Value = -20961
